# Update item names to plural forms and update bridal laces unit price
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "organza fabrics"
$ws.Range("B3").Value = "jeans fabrics"
$ws.Range("B4").Value = "suiting fabrics"
$ws.Range("B6").Value = "bridal laces"
$ws.Range("B7").Value = "black dry laces"
$ws.Range("B8").Value = "white dry laces"

# Update unit price for bridal laces
$ws.Range("C6").Value = 150

# Update the active cell selection
$ws.Range("F12").Select()
